$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = -0.3
$ws.Range("C5").Value = -8.1
$ws.Range("D5").Value = -6.5

# Replicate the cell style used for the year labels in column A (e.g. A4)
# onto the newly added A5 cell.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
